$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Date" (column B) values for rows 2-13 (re-run timestamps).
$dates = @(
    "Thu Dec 07 21:29:38 EST 2023",
    "Thu Dec 07 21:29:48 EST 2023",
    "Thu Dec 07 21:29:58 EST 2023",
    "Thu Dec 07 21:30:07 EST 2023",
    "Thu Dec 07 21:30:17 EST 2023",
    "Thu Dec 07 21:30:27 EST 2023",
    "Thu Dec 07 21:30:36 EST 2023",
    "Thu Dec 07 21:30:46 EST 2023",
    "Thu Dec 07 21:30:55 EST 2023",
    "Wed Nov 01 15:38:35 EDT 2023",
    "Wed Nov 01 15:38:45 EDT 2023",
    "Thu Dec 07 21:31:05 EST 2023"
)

# New "Execute" (column C) values for rows 2-13.
# Rows 11 & 12 (PaymentType "Extension Payments") are pulled from execution.
$executes = @("Y","Y","Y","Y","Y","Y","Y","Y","Y","DONOTRUN","DONOTRUN","Y")

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $dates[$i]
    $ws.Cells.Item($row, 3).Value = $executes[$i]
}

$ws.Range("C11:C12").Select()
